# Fix DOT multiplier and remove broken Environmental multiplier
#
# The "Menu Mock" sheet's CategoryDamageMultipliers block had a row for
# "OptionEnvironmentalMultiplier" (row 21) that is being removed entirely,
# since wall/environmental collision kills actually use DamageType.Blunt
# and never hit the "Unknown" type the Environmental multiplier was wired
# to. Deleting this row shifts every following row up by one, which also
# naturally keeps the existing "OptionDOTMultiplier" row (and everything
# below it) intact and correctly positioned.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure we're editing the "Menu Mock" sheet regardless of which sheet
# happens to be active when the workbook is opened.
$ws = $wb.Worksheets.Item("Menu Mock")

# Remove the entire row for OptionEnvironmentalMultiplier (row 21) and
# shift everything below it up by one row.
$ws.Rows("21").Delete()
